$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute("2025-07-03 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-04 Friday", 2) | Out-Null

# Update each division problem in the table by its fixed row/column position
# (avoids ambiguity from values that coincide with other cells old/new text)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "28÷6="  # was "16÷2="
$t.Cell(1,2).Range.Text = "93÷6="  # was "90÷5="
$t.Cell(1,3).Range.Text = "19÷4="  # was "95÷8="
$t.Cell(1,4).Range.Text = "80÷3="  # was "69÷5="
$t.Cell(1,5).Range.Text = "53÷5="  # was "87÷8="
$t.Cell(5,1).Range.Text = "82÷8="  # was "53÷5="
$t.Cell(5,2).Range.Text = "52÷7="  # was "77÷7="
$t.Cell(5,3).Range.Text = "70÷9="  # was "37÷9="
$t.Cell(5,4).Range.Text = "64÷4="  # was "27÷8="
$t.Cell(5,5).Range.Text = "41÷9="  # was "21÷8="
$t.Cell(9,1).Range.Text = "40÷7="  # was "36÷3="
$t.Cell(9,2).Range.Text = "86÷3="  # was "32÷3="
$t.Cell(9,3).Range.Text = "89÷5="  # was "41÷4="
$t.Cell(9,4).Range.Text = "41÷9="  # was "75÷3="
$t.Cell(9,5).Range.Text = "57÷9="  # was "79÷8="
$t.Cell(13,1).Range.Text = "51÷6="  # was "59÷9="
$t.Cell(13,2).Range.Text = "12÷7="  # was "81÷7="
$t.Cell(13,3).Range.Text = "89÷2="  # was "20÷6="
$t.Cell(13,4).Range.Text = "77÷8="  # was "90÷9="
$t.Cell(13,5).Range.Text = "35÷6="  # was "70÷9="
$t.Cell(17,1).Range.Text = "75÷9="  # was "66÷8="
$t.Cell(17,2).Range.Text = "93÷7="  # was "63÷3="
$t.Cell(17,3).Range.Text = "23÷5="  # was "34÷8="
$t.Cell(17,4).Range.Text = "39÷6="  # was "72÷3="
$t.Cell(17,5).Range.Text = "72÷8="  # was "42÷4="
